$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 412, shifting rows 412:507 down to 413:508
$ws.Rows(412).Insert()

# Populate the newly inserted row 412 with its data
$ws.Range("A412").Value = 10
$ws.Range("B412").Value = "Vega Modelo de Temuco"
$ws.Range("C412").Value = "La Araucanía"
$ws.Range("D412").Value = 45204
$ws.Range("E412").Value = 9
$ws.Range("F412").Value = "Fruta"
$ws.Range("G412").Value = 100102
$ws.Range("H412").Value = "Cítricos"
$ws.Range("I412").Value = 100102006
$ws.Range("J412").Value = "Pomelo"
$ws.Range("K412").Value = "Start Ruby"
$ws.Range("L412").Value = "Primera"
$ws.Range("M412").Value = 350
$ws.Range("N412").Value = 14000
$ws.Range("O412").Value = 15000
$ws.Range("P412").Value = 14571
$ws.Range("Q412").Value = "$/bandeja 15 kilos granel"
$ws.Range("R412").Value = "Región de O'Higgins"
$ws.Range("S412").Value = 971
$ws.Range("T412").Value = 15
